$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "42.308.37"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "2.233.29"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.59"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.00%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0962"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.13"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.47%  "
$ws.Range("E13").Value = "  +0.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.851"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").Value = "2.229.01"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("D17").Value = "42.193.66"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("E18").Value = "  +14.32%  "
$ws.Range("E19").Value = "  +1.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +38.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "231.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.08%  "
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("E28").Value = "  +3.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +19.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0805"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.117"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.90%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.38"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.01%  "
$ws.Range("E37").Value = "  +2.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "13.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.03%  "
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "63.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.82%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.02%  "
$ws.Range("E45").Value = "  +3.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.994"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("E47").Value = "  +5.97%  "
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("E49").Value = "  +1.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.78%  "
$ws.Range("E51").Value = "  -1.76%  "
